# Applies cryptos-list refresh (price / 1h-volume updates, plus a Dai/LEO row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "57.707.16"
Set-TextValue $ws "E2" "  +2.23%  "
Set-TextValue $ws "D3" "3.032.03"
Set-TextValue $ws "E3" "  +1.95%  "
Set-TextValue $ws "D5" "511.57"
Set-TextValue $ws "E5" "  +1.86%  "
Set-TextValue $ws "D6" "139.81"
Set-TextValue $ws "E6" "  +3.67%  "
Set-TextValue $ws "E7" "  +0.03%  "
Set-TextValue $ws "D8" "0.442"
Set-TextValue $ws "E8" "  +3.34%  "
Set-TextValue $ws "D9" "7.49"
Set-TextValue $ws "E9" "  +0.70%  "
Set-TextValue $ws "E10" "  +3.46%  "
Set-TextValue $ws "E11" "  +5.35%  "
Set-TextValue $ws "D12" "3.549.08"
Set-TextValue $ws "E12" "  +1.91%  "
Set-TextValue $ws "E13" "  +1.88%  "
Set-TextValue $ws "D14" "26.74"
Set-TextValue $ws "E14" "  +5.25%  "
Set-TextValue $ws "E15" "  +9.80%  "
Set-TextValue $ws "D16" "57.713.72"
Set-TextValue $ws "E16" "  +2.30%  "
Set-TextValue $ws "D17" "6.28"
Set-TextValue $ws "E17" "  +9.31%  "
Set-TextValue $ws "D18" "3.029.32"
Set-TextValue $ws "E18" "  +1.85%  "
Set-TextValue $ws "D19" "12.90"
Set-TextValue $ws "E19" "  +4.85%  "
Set-TextValue $ws "D20" "8.01"
Set-TextValue $ws "E20" "  +3.81%  "
Set-TextValue $ws "D21" "332.38"
Set-TextValue $ws "E21" "  +3.41%  "
Set-TextValue $ws "B22" "LEO"
Set-TextValue $ws "C22" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D22" "5.81"
Set-TextValue $ws "E22" "  +1.96%  "
Set-TextValue $ws "B23" "Dai"
Set-TextValue $ws "C23" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D23" "0.998"
Set-TextValue $ws "E23" "  -0.17%  "
Set-TextValue $ws "E24" "  +6.98%  "
Set-TextValue $ws "D25" "64.56"
Set-TextValue $ws "E25" "  +4.74%  "
Set-TextValue $ws "D26" "0.169"
Set-TextValue $ws "E26" "  +4.87%  "
Set-TextValue $ws "E27" "  -0.35%  "
Set-TextValue $ws "D28" "0.0₃0934"
Set-TextValue $ws "E28" "  +5.23%  "
Set-TextValue $ws "E29" "  +7.15%  "
Set-TextValue $ws "E30" "  +11.27%  "
Set-TextValue $ws "E31" "  +3.83%  "
Set-TextValue $ws "D33" "20.81"
Set-TextValue $ws "E33" "  +2.40%  "
Set-TextValue $ws "E34" "  +6.84%  "
Set-TextValue $ws "D35" "154.91"
Set-TextValue $ws "E35" "  -2.03%  "
Set-TextValue $ws "D36" "5.88"
Set-TextValue $ws "E36" "  +6.71%  "
Set-TextValue $ws "E37" "  +2.48%  "
Set-TextValue $ws "D38" "24.89"
Set-TextValue $ws "E38" "  +7.88%  "
Set-TextValue $ws "E39" "  +2.56%  "
Set-TextValue $ws "D40" "3.066.95"
Set-TextValue $ws "E40" "  +1.98%  "
Set-TextValue $ws "D41" "37.43"
Set-TextValue $ws "E41" "  +3.43%  "
Set-TextValue $ws "E42" "  +9.33%  "
Set-TextValue $ws "E43" "  +0.06%  "
Set-TextValue $ws "D44" "2.309.04"
Set-TextValue $ws "E44" "  +3.12%  "
Set-TextValue $ws "E45" "  +2.74%  "
Set-TextValue $ws "D46" "1.42"
Set-TextValue $ws "E46" "  +2.78%  "
Set-TextValue $ws "D47" "0.995"
Set-TextValue $ws "E47" "  +1.81%  "
Set-TextValue $ws "D48" "6.04"
Set-TextValue $ws "E48" "  +5.72%  "
Set-TextValue $ws "E49" "  +2.80%  "
Set-TextValue $ws "D50" "19.78"
Set-TextValue $ws "E50" "  +4.88%  "
Set-TextValue $ws "E51" "  -3.61%  "
